$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "qqwqw"
$ws.Range("E6").Value = "ffr"
$ws.Range("E6").Select() | Out-Null
